$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.1979736666666666
$ws.Range("H2").Value = 0.5939209999999999
$ws.Range("I2").Value = 0.001485174276250702
$ws.Range("J2").Value = 0.001485174276250702
$ws.Range("M2").Value = 0.1419263333333333
$ws.Range("N2").Value = 0.425779
$ws.Range("O2").Value = 0.002583058778296354
$ws.Range("P2").Value = 0.002583058778296354
$ws.Range("Q2").Value = 0.02809767660655555
$ws.Range("R2").Value = 0.252879089459
$ws.Range("S2").Value = 0.00000383629245156931
$ws.Range("T2").Value = 0.00000383629245156931
$ws.Range("G3").Value = 0.1979736666666666
$ws.Range("H3").Value = 0.5939209999999999
$ws.Range("I3").Value = 0.001485174276250702
$ws.Range("J3").Value = 0.001485174276250702
$ws.Range("O3").Value = 0.001399682868699959
$ws.Range("P3").Value = 0.001399682868699959
$ws.Range("Q3").Value = 0.01522529681744444
$ws.Range("R3").Value = 0.137027671357
$ws.Range("S3").Value = 0.000002078772991501968
$ws.Range("T3").Value = 0.000002078772991501968
$ws.Range("G4").Value = 0.1979736666666666
$ws.Range("H4").Value = 0.5939209999999999
$ws.Range("I4").Value = 0.001485174276250702
$ws.Range("J4").Value = 0.001485174276250702
$ws.Range("M4").Value = 2.613991
$ws.Range("N4").Value = 7.841973
$ws.Range("O4").Value = 0.04757462720522382
$ws.Range("P4").Value = 0.04757462720522382
$ws.Range("Q4").Value = 0.5175013829036665
$ws.Range("R4").Value = 4.657512446133
$ws.Range("S4").Value = 0.00007065661252741525
$ws.Range("T4").Value = 0.00007065661252741525
$ws.Range("G5").Value = 0.1979736666666666
$ws.Range("H5").Value = 0.5939209999999999
$ws.Range("I5").Value = 0.001485174276250702
$ws.Range("J5").Value = 0.001485174276250702
$ws.Range("M5").Value = 52.11224233333333
$ws.Range("N5").Value = 156.336727
$ws.Range("O5").Value = 0.9484426311477799
$ws.Range("P5").Value = 0.9484426311477798
$ws.Range("Q5").Value = 10.31685169295189
$ws.Range("R5").Value = 92.85166523656699
$ws.Range("S5").Value = 0.001408602598280216
$ws.Range("T5").Value = 0.001408602598280216
$ws.Range("I6").Value = 0.002310195799763575
$ws.Range("J6").Value = 0.002310195799763576
$ws.Range("M6").Value = 0.1419263333333333
$ws.Range("N6").Value = 0.425779
$ws.Range("O6").Value = 0.002583058778296354
$ws.Range("P6").Value = 0.002583058778296354
$ws.Range("Q6").Value = 0.04370607242366667
$ws.Range("R6").Value = 0.393354651813
$ws.Range("S6").Value = 0.00000596737154016267
$ws.Range("T6").Value = 0.00000596737154016267
$ws.Range("I7").Value = 0.002310195799763575
$ws.Range("J7").Value = 0.002310195799763576
$ws.Range("O7").Value = 0.001399682868699959
$ws.Range("P7").Value = 0.001399682868699959
$ws.Range("S7").Value = 0.000003233541484271677
$ws.Range("T7").Value = 0.000003233541484271678
$ws.Range("I8").Value = 0.002310195799763575
$ws.Range("J8").Value = 0.002310195799763576
$ws.Range("M8").Value = 2.613991
$ws.Range("N8").Value = 7.841973
$ws.Range("O8").Value = 0.04757462720522382
$ws.Range("P8").Value = 0.04757462720522382
$ws.Range("Q8").Value = 0.804975914459
$ws.Range("R8").Value = 7.244783230131
$ws.Range("S8").Value = 0.000109906703944826
$ws.Range("T8").Value = 0.000109906703944826
$ws.Range("I9").Value = 0.002310195799763575
$ws.Range("J9").Value = 0.002310195799763576
$ws.Range("M9").Value = 52.11224233333333
$ws.Range("N9").Value = 156.336727
$ws.Range("O9").Value = 0.9484426311477799
$ws.Range("P9").Value = 0.9484426311477798
$ws.Range("Q9").Value = 16.04791291430767
$ws.Range("R9").Value = 144.431216228769
$ws.Range("S9").Value = 0.002191088182794315
$ws.Range("T9").Value = 0.002191088182794315
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.1138773333333333
$ws.Range("H10").Value = 0.341632
$ws.Range("I10").Value = 0.0008542938511082786
$ws.Range("J10").Value = 0.0008542938511082786
$ws.Range("M10").Value = 0.1419263333333333
$ws.Range("N10").Value = 0.425779
$ws.Range("O10").Value = 0.002583058778296354
$ws.Range("P10").Value = 0.002583058778296354
$ws.Range("Q10").Value = 0.01616219236977778
$ws.Range("R10").Value = 0.145459731328
$ws.Range("S10").Value = 0.000002206691231349837
$ws.Range("T10").Value = 0.000002206691231349837
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 0.6666666666666666
$ws.Range("G11").Value = 0.1138773333333333
$ws.Range("H11").Value = 0.341632
$ws.Range("I11").Value = 0.0008542938511082786
$ws.Range("J11").Value = 0.0008542938511082786
$ws.Range("O11").Value = 0.001399682868699959
$ws.Range("P11").Value = 0.001399682868699959
$ws.Range("Q11").Value = 0.008757812238222222
$ws.Range("R11").Value = 0.078820310144
$ws.Range("S11").Value = 0.000001195740468231971
$ws.Range("T11").Value = 0.000001195740468231971
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 0.6666666666666666
$ws.Range("G12").Value = 0.1138773333333333
$ws.Range("H12").Value = 0.341632
$ws.Range("I12").Value = 0.0008542938511082786
$ws.Range("J12").Value = 0.0008542938511082786
$ws.Range("M12").Value = 2.613991
$ws.Range("N12").Value = 7.841973
$ws.Range("O12").Value = 0.04757462720522382
$ws.Range("P12").Value = 0.04757462720522382
$ws.Range("Q12").Value = 0.2976743244373333
$ws.Range("R12").Value = 2.679068919936
$ws.Range("S12").Value = 0.00004064271149019134
$ws.Range("T12").Value = 0.00004064271149019134
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 0.6666666666666666
$ws.Range("G13").Value = 0.1138773333333333
$ws.Range("H13").Value = 0.341632
$ws.Range("I13").Value = 0.0008542938511082786
$ws.Range("J13").Value = 0.0008542938511082786
$ws.Range("M13").Value = 52.11224233333333
$ws.Range("N13").Value = 156.336727
$ws.Range("O13").Value = 0.9484426311477799
$ws.Range("P13").Value = 0.9484426311477798
$ws.Range("Q13").Value = 5.934403190940444
$ws.Range("R13").Value = 53.409628718464
$ws.Range("S13").Value = 0.0008102487079185055
$ws.Range("T13").Value = 0.0008102487079185053
$ws.Range("G14").Value = 132.6801566666667
$ws.Range("H14").Value = 398.04047
$ws.Range("I14").Value = 0.9953503360728775
$ws.Range("J14").Value = 0.9953503360728775
$ws.Range("M14").Value = 0.1419263333333333
$ws.Range("N14").Value = 0.425779
$ws.Range("O14").Value = 0.002583058778296354
$ws.Range("P14").Value = 0.002583058778296354
$ws.Range("Q14").Value = 18.83080814179223
$ws.Range("R14").Value = 169.47727327613
$ws.Range("S14").Value = 0.002571048423073272
$ws.Range("T14").Value = 0.002571048423073272
$ws.Range("G15").Value = 132.6801566666667
$ws.Range("H15").Value = 398.04047
$ws.Range("I15").Value = 0.9953503360728775
$ws.Range("J15").Value = 0.9953503360728775
$ws.Range("O15").Value = 0.001399682868699959
$ws.Range("P15").Value = 0.001399682868699959
$ws.Range("Q15").Value = 10.20385590188778
$ws.Range("R15").Value = 91.83470311699001
$ws.Range("S15").Value = 0.001393174813755953
$ws.Range("T15").Value = 0.001393174813755953
$ws.Range("G16").Value = 132.6801566666667
$ws.Range("H16").Value = 398.04047
$ws.Range("I16").Value = 0.9953503360728775
$ws.Range("J16").Value = 0.9953503360728775
$ws.Range("M16").Value = 2.613991
$ws.Range("N16").Value = 7.841973
$ws.Range("O16").Value = 0.04757462720522382
$ws.Range("P16").Value = 0.04757462720522382
$ws.Range("Q16").Value = 346.8247354052567
$ws.Range("R16").Value = 3121.42261864731
$ws.Range("S16").Value = 0.04735342117726139
$ws.Range("T16").Value = 0.04735342117726139
$ws.Range("G17").Value = 132.6801566666667
$ws.Range("H17").Value = 398.04047
$ws.Range("I17").Value = 0.9953503360728775
$ws.Range("J17").Value = 0.9953503360728775
$ws.Range("M17").Value = 52.11224233333333
$ws.Range("N17").Value = 156.336727
$ws.Range("O17").Value = 0.9484426311477799
$ws.Range("P17").Value = 0.9484426311477798
$ws.Range("Q17").Value = 6914.260477037966
$ws.Range("R17").Value = 62228.3442933417
$ws.Range("S17").Value = 0.9440326916587869
$ws.Range("T17").Value = 0.9440326916587868

Write-Output "Updated 182 cells"
